# Updated Block Diagram with Ownership. --Spencer.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> Points helper (1 pt = 12700 EMU)
function EMU($v) { return $v / 12700 }

# ---------------------------------------------------------------------
# 1) Split "Monomial ordering from partition info" into three runs so
#    the word "Grobner" can carry its own (spell-check flagged) run.
# ---------------------------------------------------------------------
$ordShape = $s.Shapes.Item(5)   # "Rectangle 9" / Monomial ordering...
$tr = $ordShape.TextFrame.TextRange
$tr.Text = "Monomial ordering/"
$r2 = $tr.InsertAfter("Grobner")
$r3 = $r2.InsertAfter(" from partition info")

# ---------------------------------------------------------------------
# 2) Add the three ownership labels as borderless rectangles, grouping
#    the existing workflow boxes under each contributor's name.
# ---------------------------------------------------------------------
$template = $s.Shapes.Item(1)

# --- Drew (Perl) -------------------------------------------------------
$r1 = $template.Duplicate().Item(1)
$r1.Name = "Rectangle 22"
$r1.Left = EMU(304800)
$r1.Top = EMU(152400)
$r1.Width = EMU(7543800)
$r1.Height = EMU(2209800)
$r1.Fill.Visible = 0
$r1.TextFrame.VerticalAnchor = 4    # msoAnchorBottom
$r1.TextFrame.HorizontalAnchor = 0  # msoAnchorNone (anchorCtr="0")
$r1.TextFrame.TextRange.Text = "Drew (Perl)"
$r1.TextFrame.TextRange.Font.Color.ObjectThemeColor = 1  # tx1/dk1

# --- Matt (Singular) ----------------------------------------------------
$r2b = $template.Duplicate().Item(1)
$r2b.Name = "Rectangle 23"
$r2b.Left = EMU(228600)
$r2b.Top = EMU(3299388)
$r2b.Width = EMU(2362200)
$r2b.Height = EMU(2209801)
$r2b.Fill.Visible = 0
$r2b.TextFrame.VerticalAnchor = 4
$r2b.TextFrame.HorizontalAnchor = 0
$r2b.TextFrame.TextRange.Text = "Matt (Singular)"
$r2b.TextFrame.TextRange.Font.Color.ObjectThemeColor = 1

# --- Spencer(Perl) -------------------------------------------------------
$r3b = $template.Duplicate().Item(1)
$r3b.Name = "Rectangle 24"
$r3b.Left = EMU(2637446)
$r3b.Top = EMU(3299389)
$r3b.Width = EMU(3915754)
$r3b.Height = EMU(2209800)
$r3b.Fill.Visible = 0
$r3b.TextFrame.VerticalAnchor = 4
$r3b.TextFrame.HorizontalAnchor = 0
$r3b.TextFrame.TextRange.Text = "Spencer(Perl)"
$r3b.TextFrame.TextRange.Font.Color.ObjectThemeColor = 1
